# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 28, shifting the existing
# rows 28:33 down to 29:34 (dimension grows from A1:T33 to A1:T34).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing data rows (old 28:33) down by one row.
$ws.Rows("28:28").Insert()

# Populate the newly inserted row 28 with this week's record. The
# descriptive / categorical columns repeat the same market & product
# info as the rest of the block; only the date, volume and the price
# columns change week to week.
$ws.Range("A28").Value = 8
$ws.Range("B28").Value = "Terminal La Palmera de La Serena"
$ws.Range("C28").Value = "Coquimbo"
$ws.Range("D28").Value = 44782
$ws.Range("E28").Value = 4
$ws.Range("F28").Value = "Fruta"
$ws.Range("G28").Value = 100108
$ws.Range("H28").Value = "Tropicales y subtropicales"
$ws.Range("I28").Value = 100108007
$ws.Range("J28").Value = "Coco"
$ws.Range("K28").Value = "Sin especificar"
$ws.Range("L28").Value = "Primera"
$ws.Range("M28").Value = 200
$ws.Range("N28").Value = 23500
$ws.Range("O28").Value = 24000
$ws.Range("P28").Value = 23750
$ws.Range("Q28").Value = "$/malla 20 unidades"
$ws.Range("R28").Value = "Perú"
$ws.Range("S28").Value = 1188
$ws.Range("T28").Value = 20
